# Insert a new data row at row 268 (pushing the existing rows 268-361 down
# to 269-362) and populate the new row with the latest price-report entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(268).Insert()

$ws.Range("A268").Value = 10
$ws.Range("B268").Value = "Vega Modelo de Temuco"
$ws.Range("C268").Value = "La Araucanía"
$ws.Range("D268").Value = 44809
$ws.Range("E268").Value = 9
$ws.Range("F268").Value = 100112017
$ws.Range("G268").Value = "Apio"
$ws.Range("H268").Value = "Americana (o)"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 200
$ws.Range("K268").Value = 12000
$ws.Range("L268").Value = 12000
$ws.Range("M268").Value = 12000
$ws.Range("N268").Value = "$/docena de matas"
$ws.Range("O268").Value = "Provincia del Elquí"
$ws.Range("P268").Value = 2000
$ws.Range("Q268").Value = 6
$ws.Range("R268").Value = "Hortaliza"
